$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# C3 loses its "credentials" wrapper - it becomes the canonical JSON
# blob that C4 and C5 now reference via formula (generation with
# references), instead of each row carrying its own hard-coded copy.
# -----------------------------------------------------------------
$json = @"
 {  "countryCode": "280",
    "custID": "",
    "custID2": "",
    "hbciVersion": "0",
    "language": 1,
    "url": "bawagPSK.js",
    "userID": "64769092",
    "pin": "38567" }
"@

$ws.Range("C3").Value2 = $json
$ws.Range("C4").Formula = "=C3"
$ws.Range("C5").Formula = "=C4"

# -----------------------------------------------------------------
# Row 3: D3/E3 now illustrate the curly-quote convention, F3/G3 are
# new columns repeating the "||" separator + generated JSON.
# -----------------------------------------------------------------
$ws.Range("D3").Value2 = "$([char]0x201C)38567$([char]0x201D)"
$ws.Range("E3").Value2 = "$([char]0x201D)12345$([char]0x201D)"
$ws.Range("F3").Value2 = "||"
$ws.Range("G3").Value2 = $json

# Row 4: D4/E4 are cleared (the old "|" separator + duplicate JSON
# move out to the new F4/G4 columns).
$ws.Range("D4").Value2 = ""
$ws.Range("E4").Value2 = ""
$ws.Range("F4").Value2 = "|"
$ws.Range("G4").Value2 = $json

# Row 5: D5 reuses the curly-quote convention string, E5 is cleared,
# F5/G5 are the new "||" separator + generated JSON columns.
$ws.Range("D5").Value2 = "$([char]0x201C)38567$([char]0x201D)"
$ws.Range("E5").Value2 = ""
$ws.Range("F5").Value2 = "||"
$ws.Range("G5").Value2 = $json

# -----------------------------------------------------------------
# View state: scroll back to the top-left and select E3.
# -----------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.TabRatio = 984
$ws.Range("E3").Select()
